$d = $word.ActiveDocument

# 1) The sentence "A little annoyed, I consider not opening the door, but
#    after weighing the possible consequences in my mind I end up opening
#    it anyways." was previously split across several <w:r> runs. Running a
#    plain-text Find & Replace over the whole sentence collapses it back
#    down into a single run (Word rewrites the matched range as one run).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found1 = $find.Execute(
    "A little annoyed, I consider not opening the door, but after weighing the possible consequences in my mind I end up opening it anyways.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A little annoyed, I consider not opening the door, but after weighing the possible consequences in my mind I end up opening it anyways.",
    2
)

# 2) Petra's expression changes from "raised_eyebrow" to "expressionless".
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$found2 = $find2.Execute(
    "Petra (neutral raised_eyebrow): You do, huh?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Petra (neutral expressionless): You do, huh?",
    2
)
